$d = $word.ActiveDocument

# 1. Insert a new paragraph at the very beginning with red text "Daniel’s edits in red"
$startRng = $d.Range(0, 0)
$startRng.InsertBefore("Daniel’s edits in red`r")
$firstPara = $d.Paragraphs(1).Range
$firstPara.Font.Color = 255

# 2. Merge the two runs "Incremental Development" + " Process model diagram"
#    into a single run containing the combined text.
$d.Content.Find.Execute("Incremental Development Process model diagram", $true, $false, $false, $false, $false, $true, 1, $false, "Incremental Development Process model diagram", 2)

# 3. Mark the run containing the drawing as noProof
$shp = $d.InlineShapes(1)
$shp.Range.Font.NoProofing = 1

# 4. Add a new red list item after "Adjust and fix functions based on feedback from beta"
$targetIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.TrimEnd("`r") -eq "Adjust and fix functions based on feedback from beta") {
        $targetIndex = $i
    }
}
$targetPara = $d.Paragraphs($targetIndex)
$endRng = $targetPara.Range
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
$newPara = $d.Paragraphs($targetIndex + 1)
$newPara.Range.ListFormat.ListIndent()
$textRng = $newPara.Range
$textRng.Collapse(0)
$textRng.MoveEnd(1, -1)
$textRng.Text = "Also check for user access adjustments, making sure that everyone that needs access has it"
$textRng.Font.Color = 255
